# Update the "blast furnaces and steel mills" data sheet with corrected/expanded
# BLS figures:
#  - the women's employment figure for 1958 (C2) is not available -> "n/a"
#  - every year (1958-1974, rows 2-18) now has a complete row of the 12 series
#  - the stale, ragged leftover rows (19-32) from the previous partial layout
#    are removed now that the data is densely packed into rows 2-18
#  - leave the selection on D2, matching the author's final view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Year, B all_employee, C women, D prod, E prod_wk_earn, F prod_hr_earn,
#          G prod_wk_hr, H prod_wk_ot, I acc, J newhire, K sep, L quit, M layoff
$rows = @(
  @(1958, 531.4, "n/a", 432.2, 108.54, 2.91, 37.3, 0.8, 2.8, 0.3, 3.2, 0.2, 2.6),
  @(1959, 515.3, 20.6, 414.7, 123.38, 3.1, 39.8, 2.1, 3.3, 1.6, 1.4, 0.5, 0.4),
  @(1960, 577.1, 22.2, 470.5, 117.04, 3.08, 38.0, 1.2, 2.0, 0.3, 4.3, 0.3, 3.5),
  @(1961, 526.5, 20.4, 424.7, 123.84, 3.2, 38.7, 1.2, 3.5, 0.4, 2.4, 0.3, 1.5),
  @(1962, 522.3, 19.9, 421.4, 128.31, 3.29, 39.0, 1.3, 2.6, 0.4, 3.7, 0.3, 2.8),
  @(1963, 520.0, 18.8, 424.6, 134.4, 3.36, 40.0, 1.8, 3.2, 0.7, 2.8, 0.4, 1.7),
  @(1964, 556.7, 19.1, 458.4, 140.15, 3.41, 41.1, 2.3, 2.9, 1.5, 1.8, 0.6, 0.5),
  @(1965, 580.2, 19.6, 477.4, 141.86, 3.46, 41.0, 2.5, 2.3, 1.3, 3.0, 0.9, 1.3),
  @(1966, 571.4, 19.9, 467.2, 145.71, 3.58, 40.7, 2.4, 2.9, 1.7, 2.4, 1.1, 0.5),
  @(1967, 555.5, 20.5, 447.8, 145.16, 3.62, 40.1, 2.0, 2.5, 1.2, 2.5, 0.8, 0.9),
  @(1968, 555.5, 20.7, 445.1, 155.86, 3.82, 40.8, 2.8, 3.0, 1.5, 3.5, 1.1, 1.4),
  @(1969, 561.1, 21.5, 450.0, 168.51, 4.09, 41.2, 3.0, 3.3, 2.2, 2.7, 1.3, 0.3),
  @(1970, 546.3, 21.9, 437.1, 168.38, 4.22, 39.9, 2.2, 2.7, 1.4, 3.3, 1.1, 1.2),
  @(1971, 497.3, 20.1, 394.9, 181.43, 4.57, 39.7, 2.2, 3.5, 1.0, 4.6, 0.7, 3.0),
  @(1972, 491.9, 18.8, 393.6, 210.12, 5.15, 40.8, 2.4, 3.1, 0.9, 2.2, 0.6, 0.8),
  @(1973, 518.4, 20.3, 418.1, 230.74, 5.56, 41.5, 3.2, 2.5, 1.7, 2.1, 0.9, 0.4),
  @(1974, 522.0, 22.1, 418.7, 263.08, 6.37, 41.3, 3.2, 2.0, 1.1, 2.3, 0.7, 0.7)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $rowArr = $rows[$i]
  $r = $i + 2
  for ($j = 0; $j -lt $rowArr.Count; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $rowArr[$j]
  }
}

# drop the old leftover rows below the (now fully populated) data block
$ws.Range("A19:M32").Clear()

$ws.Range("D2").Select()
